$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''30.176.26'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  -3.22%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''1.863.05'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  -3.97%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('D4').Value = '''1.000'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '''  -0.01%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''234.06'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  -3.37%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''1.000'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  -0.02%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = '''0.4669'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '''  -2.68%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = '''0.2829'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '''  -2.90%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = '''0.06549'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '''  -3.51%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = '''20.19'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  +0.08%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = '''0.07813'
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Value = '''96.02'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''  -7.95%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = '''1.856.26'
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Value = '''5.124'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '''  -3.41%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''0.6703'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  -4.09%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = '''280.57'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''  -5.38%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''30.211.26'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  -3.07%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''1.0000'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  -0.02%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = '''5.449'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '''  -2.26%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''12.65'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  -2.78%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = '''2.099.85'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '''  -4.65%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''0.000007254'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  -4.64%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = '''0.9994'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '''  -0.10%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = '''6.147'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '''  -4.25%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''9.319'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  -2.57%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = '''165.26'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '''  -2.25%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('E27').Value = '''  -4.65%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('E28').Value = '''  -9.12%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('E29').Value = '''  -3.49%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = '''0.09604'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '''  -4.78%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = '''4.408'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '''  -4.67%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''1.468'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  -4.39%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('E33').Value = '''  -5.15%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''0.04665'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  -3.76%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = '''0.7016'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '''  -5.05%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('E36').Value = '''  -3.10%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = '''2.709'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '''  -0.53%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = '''0.01867'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '''  -4.97%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = '''6.278'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '''  -8.01%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = '''2.527'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  -4.05%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').Value = '''72.53'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '''  -5.07%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = '''0.8528'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''  -2.34%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('E43').Value = '''  -5.49%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = '''0.4164'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''  -4.72%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('E45').Value = '''  -0.04%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = '''103.17'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '''  -2.49%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = '''985.28'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '''  -4.66%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = '''7.137'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '''  -5.98%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = '''9.191'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '''  -0.73%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('D51').Value = '''0.1138'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '''  -5.88%  '
$ws.Range('E51').Style = "Normal"
